$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header
$ws.Range("B1").Value = "mapsto"

# Update column B labels (better stage names)
$ws.Range("B2").Value = "wake"
$ws.Range("B3").Value = "rem"
$ws.Range("B4").Value = "stage1"
$ws.Range("B5").Value = "stage2"
$ws.Range("B6").Value = "sws"
$ws.Range("B7").Value = "movement"
$ws.Range("B8").Value = "unknown"
$ws.Range("B9").Value = "artifact"

# Remove the now-unused "meaning" column C entirely
$ws.Columns.Item(3).Delete()

# Update selection to match target state
$ws.Range("D3").Select()
